# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Mon Nov 13 10:40:08 UTC 2023 with GitHub Actions")
#
# All Price (D) / Volume(1h) (E) cells on this sheet are stored as plain TEXT
# (t="inlineStr" in the OOXML, e.g. "246.13" or "  -1.43%  "), not numbers. For any
# new value that looks like a plain decimal number (a single "." and digits, e.g.
# "246.13"), Excel would otherwise auto-convert the cell to a numeric value when the
# .Value is assigned, so we force those specific cells to Text format ("@") first so
# the write round-trips as a string, matching the source data untouched otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.970.50'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.051.28'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.13'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.95'
$ws.Range('E7').Value = '  -3.71%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0781'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.31'
$ws.Range('E12').Value = '  -4.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.878'
$ws.Range('E13').Value = '  +4.97%  '
$ws.Range('D14').Value = '2.351.12'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.67'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').Value = '2.090.43'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '17.96'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '36.936.16'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.56'
$ws.Range('E19').Value = '  -3.63%  '
$ws.Range('D20').Value = '0.0₃0894'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.42'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.17'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.36'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.12'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.03'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.55'
$ws.Range('E29').Value = '  +14.77%  '
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.11'
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.70'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.34'
$ws.Range('E34').Value = '  +3.96%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.84'
$ws.Range('E36').Value = '  +5.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0810'
$ws.Range('E37').Value = '  -8.38%  '
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.16'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('E43').Value = '  -11.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.14'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.92'
$ws.Range('E45').Value = '  -5.20%  '
$ws.Range('D46').Value = '1.305.22'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  -5.93%  '
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.75'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('D50').Value = '2.236.04'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.65'
$ws.Range('E51').Value = '  +1.32%  '
